# Tutorial 6 solution update: switch date separators from "/" to "-"
# and correct the attendance tallies for the first (28-07-2022) and
# sixth (01-09-2022) session rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> corrected date string (DD-MM-YYYY)
$dates = @{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

# Force Text format first so the COM layer doesn't reinterpret these
# day-month-year strings as date serials (and potentially swap day/month).
foreach ($row in $dates.Keys) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$row]
}

# Row 3 (28-07-2022): mark the attendance as real (D) and invalid (G) instead of all zero.
$ws.Cells.Item(3, 4).Value = 1   # D3 Total Attendance Count
$ws.Cells.Item(3, 7).Value = 1   # G3 Invalid

# Row 13 (01-09-2022): record real attendance (D, E) and clear the absence flag (H).
$ws.Cells.Item(13, 4).Value = 1  # D13 Total Attendance Count
$ws.Cells.Item(13, 5).Value = 1  # E13 Real
$ws.Cells.Item(13, 8).Value = 0  # H13 Absent
